$d = $word.ActiveDocument

# Locate the target paragraphs by their leading text rather than a hard-coded
# index, so the script is resilient to any incidental paragraph reshuffling.
$rilevazioneIndex = 0
$replicaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($rilevazioneIndex -eq 0 -and $t.StartsWith("Rilevazione")) {
        $rilevazioneIndex = $i
    }
    if ($replicaIndex -eq 0 -and $t.StartsWith("Replica ")) {
        $replicaIndex = $i
    }
}
if ($rilevazioneIndex -eq 0) { $rilevazioneIndex = 6 }
if ($replicaIndex -eq 0) { $replicaIndex = 9 }

# --- Paragraph "Rilevazione ( ... )": split DataRil/RespRil runs to add "ev",
# and append ", ClasseRilev(CLASSE), IndividuoRilev(PERSONA)" before the closing ")" ---
$rilevazione = $d.Paragraphs.Item($rilevazioneIndex).Range
$rilevazioneXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="00000006" w14:textId="5D174EC9" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:right="-182" w:hanging="285"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Rilevazione</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodR</w:t></w:r><w:r><w:t>, RespIns</w:t></w:r><w:r w:rsidR="002B675B"><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>0</w:t></w:r><w:r><w:t>, DataRil</w:t></w:r><w:r><w:t>ev</w:t></w:r><w:r><w:t>, DataIns, RespRil</w:t></w:r><w:r><w:t>ev</w:t></w:r><w:r><w:t>, ModAcquisizione, InfoAmb</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>INFOAMBIENTALI</w:t></w:r><w:r w:rsidR="002614FF"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="002614FF" w:rsidRPr="002614FF"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Dispositivo</w:t></w:r><w:r w:rsidR="002614FF" w:rsidRPr="002614FF"><w:rPr><w:color w:val="FF0000"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>DISPOSITIVO</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>, ClasseRilev</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>CLASSE</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>, IndividuoRilev</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>PERSONA</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
'@
$rilevazione.InsertXML($rilevazioneXml)

# --- Paragraph "Replica ( ... )": merge the split "O"/"RTO" superscript runs into "ORTO" ---
$replica = $d.Paragraphs.Item($replicaIndex).Range
$replicaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="00000009" w14:textId="16228274" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Replica</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="008E02E9"><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r><w:r w:rsidRPr="008E02E9"><w:rPr><w:u w:val="single"/></w:rPr><w:t>Replica</w:t></w:r><w:r w:rsidR="00181598"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> Gruppo, DataDimora, Esposizione</w:t></w:r><w:r w:rsidR="008E02E9" w:rsidRPr="008E02E9"><w:t>,</w:t></w:r><w:r w:rsidR="002639AD"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="002639AD"><w:t>Specie</w:t></w:r><w:r w:rsidR="008E02E9" w:rsidRPr="008E02E9"><w:t>Pianta</w:t></w:r><w:r w:rsidR="008E02E9" w:rsidRPr="008E02E9"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>SPECIE</w:t></w:r><w:r w:rsidR="006143A5"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="002639AD"><w:t>C</w:t></w:r><w:r w:rsidR="006143A5"><w:t>lasseDimora</w:t></w:r><w:r w:rsidR="006143A5"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>CLASS</w:t></w:r><w:r w:rsidR="006143A5" w:rsidRPr="008E02E9"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>E</w:t></w:r><w:r w:rsidR="002614FF"><w:t>,</w:t></w:r><w:r w:rsidR="002614FF" w:rsidRPr="002614FF"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="002614FF"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Orto</w:t></w:r><w:r w:rsidR="002614FF"><w:rPr><w:color w:val="FF0000"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>ORTO</w:t></w:r><w:r w:rsidR="00E30FA7"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00E30FA7" w:rsidRPr="002614FF"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Dispositivo</w:t></w:r><w:r w:rsidR="00E30FA7" w:rsidRPr="002614FF"><w:rPr><w:color w:val="FF0000"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>DISPOSITIVO</w:t></w:r><w:r w:rsidR="006143A5"><w:t>)</w:t></w:r></w:p>
'@
$replica.InsertXML($replicaXml)
